# 3.3.1.1 --> 3.4.0 update hydgn/HPEbP
#
# Substantive change: on the "HPEbP" worksheet, the natural-gas reforming
# efficiency formula in B3 drops the "+46" term, going from
# 118/(162+2+46) to 118/(162+2). Every other cell in that row (C3:AI3)
# simply carries the value forward (="previous cell" / shared formula),
# so they recompute automatically once B3 changes - no need to touch them
# directly.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsIEA   = $wb.Worksheets.Item("IEA Data")
$wsHPEbP = $wb.Worksheets.Item("HPEbP")

# The actual data/formula fix.
$wsHPEbP.Range("B3").Formula = "=118/(162+2)"

# Leave the workbook's view state the way it was after the author's last
# save: HPEbP active/selected, with the selection cursors on each sheet
# where they were left.
[void]$wsAbout.Activate()
[void]$wsAbout.Range("B14").Select()

[void]$wsIEA.Activate()
[void]$wsIEA.Range("E18").Select()

[void]$wsHPEbP.Activate()
[void]$wsHPEbP.Range("K5").Select()
